$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as literal text (e.g. "96.822.59",
# "0.0000255") where "." is used as a thousands/format separator, not a true decimal
# point. Assigning such a string directly to .Value lets the COM layer auto-coerce it
# into a real number, corrupting the text. Forcing the cell to Text format ("@")
# before the assignment keeps the literal string; ClearFormats() afterwards removes
# that temporary number-format override so the cell style stays at the workbook
# default (no explicit style), matching the original unstyled cells.

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '96.822.59'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +5.53%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.608.49'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +8.84%  '

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '239.30'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.23%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '640.02'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.50%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.49'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +7.51%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.405'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.74%  '

# Row 9
$ws.Range('E9').Value = '  -0.20%  '

# Row 10
$ws.Range('E10').Value = '  +8.76%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '3.605.80'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +8.83%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '43.50'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.85%  '

# Row 13
$ws.Range('E13').Value = '  +4.51%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.41'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +7.40%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.288.56'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +8.92%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '96.716.00'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +5.61%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000255'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +5.10%  '

# Row 18
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.612.02'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +9.23%  '

# Row 19
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.19'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +22.14%  '

# Row 20
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.75'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.61%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '18.14'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +5.46%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.499'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +12.68%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '517.47'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +5.70%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.47'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.87%  '

# Row 25
$ws.Range('E25').Value = '  +10.94%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.74'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +11.07%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '97.67'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.84%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.51'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +6.42%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.12'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +20.13%  '

# Row 30
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.145'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.29%  '

# Row 31
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '11.64'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +6.79%  '

# Row 32
$ws.Range('E32').Value = '  -0.07%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.183'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.98%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.992'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.56%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '30.51'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.45%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.574'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +10.43%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '579.80'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +7.70%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '7.92'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +8.32%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.50'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +10.20%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.153'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.59%  '

# Row 41
$ws.Range('E41').Value = '  +0.02%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.928'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +8.15%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.76'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +6.46%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0434'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +6.72%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '23.81'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.39%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.71'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +7.33%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.21'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +6.48%  '

# Row 48
$ws.Range('E48').Value = '  -1.13%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '53.94'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.53%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.18'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.23%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.13'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.27%  '
